$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D) updates ---
# Use NumberFormat "@" (Text) while assigning so Excel does not
# reinterpret values like "19.58" or "0.0840" as numbers, then
# restore the original cell style so no formatting changes leak in.
$cell = $ws.Range('D2')
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '26.662.86'
$cell.Style = $savedStyle

$cell = $ws.Range('D3')
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.597.90'
$cell.Style = $savedStyle

$cell = $ws.Range('D5')
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '211.58'
$cell.Style = $savedStyle

$cell = $ws.Range('D9')
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.247'
$cell.Style = $savedStyle

$cell = $ws.Range('D10')
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '19.58'
$cell.Style = $savedStyle

$cell = $ws.Range('D11')
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0840'
$cell.Style = $savedStyle

$cell = $ws.Range('D12')
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.822.00'
$cell.Style = $savedStyle

$cell = $ws.Range('D13')
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.660.74'
$cell.Style = $savedStyle

$cell = $ws.Range('D16')
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '65.07'
$cell.Style = $savedStyle

$cell = $ws.Range('D17')
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '26.643.22'
$cell.Style = $savedStyle

$cell = $ws.Range('D20')
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '208.93'
$cell.Style = $savedStyle

$cell = $ws.Range('D21')
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '7.02'
$cell.Style = $savedStyle

$cell = $ws.Range('D22')
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '4.29'
$cell.Style = $savedStyle

$cell = $ws.Range('D23')
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.33'
$cell.Style = $savedStyle

$cell = $ws.Range('D25')
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '144.34'
$cell.Style = $savedStyle

$cell = $ws.Range('D34')
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.290.41'
$cell.Style = $savedStyle

$cell = $ws.Range('D39')
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.834'
$cell.Style = $savedStyle

$cell = $ws.Range('D40')
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.06'
$cell.Style = $savedStyle

$cell = $ws.Range('D43')
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.785'
$cell.Style = $savedStyle

$cell = $ws.Range('D44')
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '63.47'
$cell.Style = $savedStyle

$cell = $ws.Range('D45')
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '1.734.34'
$cell.Style = $savedStyle

$cell = $ws.Range('D46')
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '90.77'
$cell.Style = $savedStyle

$cell = $ws.Range('D49')
$savedStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0509'
$cell.Style = $savedStyle

# --- Volume(1h) column (E) updates ---
$ws.Range('E2').Value = '  -0.06%  '
$ws.Range('E3').Value = '  -0.08%  '
$ws.Range('E4').Value = '  +0.29%  '
$ws.Range('E5').Value = '  +0.19%  '
$ws.Range('E6').Value = '  +0.49%  '
$ws.Range('E7').Value = '  +0.24%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').Value = '  +0.79%  '
$ws.Range('E10').Value = '  -0.63%  '
$ws.Range('E11').Value = '  +0.36%  '
$ws.Range('E12').Value = '  -0.02%  '
$ws.Range('E13').Value = '  +3.89%  '
$ws.Range('E14').Value = '  -0.15%  '
$ws.Range('E15').Value = '  +0.41%  '
$ws.Range('E16').Value = '  +0.32%  '
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('E20').Value = '  -0.45%  '
$ws.Range('E21').Value = '  +3.56%  '
$ws.Range('E22').Value = '  +0.51%  '
$ws.Range('E23').Value = '  +1.66%  '
$ws.Range('E24').Value = '  +0.94%  '
$ws.Range('E25').Value = '  -1.38%  '
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('E27').Value = '  -0.80%  '
$ws.Range('E28').Value = '  -0.71%  '
$ws.Range('E29').Value = '  -0.14%  '
$ws.Range('E30').Value = '  +2.32%  '
$ws.Range('E31').Value = '  +0.23%  '
$ws.Range('E32').Value = '  +0.15%  '
$ws.Range('E33').Value = '  +1.54%  '
$ws.Range('E34').Value = '  -0.48%  '
$ws.Range('E35').Value = '  -8.18%  '
$ws.Range('E36').Value = '  +0.53%  '
$ws.Range('E37').Value = '  +0.44%  '
$ws.Range('E38').Value = '  -0.73%  '
$ws.Range('E39').Value = '  -0.98%  '
$ws.Range('E40').Value = '  +17.74%  '
$ws.Range('E41').Value = '  +2.16%  '
$ws.Range('E43').Value = '  -0.26%  '
$ws.Range('E44').Value = '  -0.78%  '
$ws.Range('E46').Value = '  +0.78%  '
$ws.Range('E47').Value = '  -3.47%  '
$ws.Range('E48').Value = '  +1.40%  '
$ws.Range('E49').Value = '  +0.96%  '
$ws.Range('E50').Value = '  +0.30%  '
$ws.Range('E51').Value = '  -1.10%  '
